# Add a "Segment / Count of patients" summary table in columns M:N,
# mirroring the Frail Elders ... Relatively healthy rows already present
# in columns A:B, with a styled (bold white-on-blue) header row and a
# thin border + centered body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values -----------------------------------------------------------
$ws.Range("M1").Value = "Segment"
$ws.Range("N1").Value = "Count of patients"

$segments = @(
    "Frail Elders",
    "People with complex multimorbidity",
    "People with major disability",
    "People with minor morbidity",
    "People with moderate morbidity",
    "People with serious mental illness",
    "Relatively healthy"
)
$counts = @(144, 241, 48, 8, 102, 91, 1)

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 13).Value = $segments[$i]
    $ws.Cells.Item($row, 14).Value = $counts[$i]
}

# ---- Formatting ---------------------------------------------------------
$headerRange = $ws.Range("M1:N1")
$bodyRange = $ws.Range("M2:N8")

# Header: bold white text on a solid blue fill, thin border, centered.
$headerRange.Font.Color = 0
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 12611584
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108

# Body: black text, thin border, centered.
$bodyRange.Font.Color = 0
$bodyRange.Borders.LineStyle = 1
$bodyRange.HorizontalAlignment = -4108

# ---- Selection ----------------------------------------------------------
$ws.Range("M1:N8").Select()
